$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old rows (2-5) entirely, leaving just the data in row 1
$ws.Rows("2:5").Delete()

# Set the new header/data values in row 1
$ws.Range("A1").Value = "asd"
$ws.Range("B1").Value = "5fa285e1bebe0a6623e33afc04a1fbd5"
